# Generate Report for Handoff
# Update the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps for the
# 94474555-9a87-4390-bad8-eca5e0787dd8.md row after a fresh handoff xliff generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 7 is 94474555-9a87-4390-bad8-eca5e0787dd8.md, column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G7").Value = "2016-09-07 02:51:48"

# zh-cn sheet: row 7 is 94474555-9a87-4390-bad8-eca5e0787dd8.md, column H = "Latest Handoff Datetime"
$wsZhCn.Range("H7").Value = "2016-09-07 02:51:43"

# de-de sheet: row 7 is 94474555-9a87-4390-bad8-eca5e0787dd8.md, column H = "Latest Handoff Datetime"
$wsDeDe.Range("H7").Value = "2016-09-07 02:51:48"
